# Update profit files after running on 2025-09-10
# Append the newly computed allocation row (row 9) to the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use an apostrophe-prefixed value so the date-like string "09/10/2025" is
# entered as literal text (matching the existing Date column entries)
# instead of being auto-converted into a date serial number, then clear
# the resulting "text" formatting so the cell keeps the sheet's default
# (unstyled) look, just like the other rows in column A.
$ws.Range("A9").Value = "'09/10/2025"
$ws.Range("A9").ClearFormats()

$ws.Range("B9").Value = 0.1250185577428985
$ws.Range("C9").Value = 0.8749814422571015
